$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf18"
$ws.Cells.Item(2,3).Value = "Fgfr3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.1601763333333333
$ws.Cells.Item(2,8).Value = 0.480529
$ws.Cells.Item(2,9).Value = 0.01412814675921196
$ws.Cells.Item(2,10).Value = 0.01412814675921196
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 5.751166666666666
$ws.Cells.Item(2,14).Value = 17.2535
$ws.Cells.Item(2,15).Value = 0.7405222614421495
$ws.Cells.Item(2,16).Value = 0.7405222614421495
$ws.Cells.Item(2,17).Value = 0.9212007890555555
$ws.Cells.Item(2,18).Value = 8.290807101499999
$ws.Cells.Item(2,19).Value = 0.01046220718811821
$ws.Cells.Item(2,20).Value = 0.01046220718811821
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf18"
$ws.Cells.Item(3,3).Value = "Fgfr3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.1601763333333333
$ws.Cells.Item(3,8).Value = 0.480529
$ws.Cells.Item(3,9).Value = 0.01412814675921196
$ws.Cells.Item(3,10).Value = 0.01412814675921196
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.5698483333333333
$ws.Cells.Item(3,14).Value = 1.709545
$ws.Cells.Item(3,15).Value = 0.07337387367415998
$ws.Cells.Item(3,16).Value = 0.07337387367416
$ws.Cells.Item(3,17).Value = 0.09127621658944444
$ws.Cells.Item(3,18).Value = 0.821485949305
$ws.Cells.Item(3,19).Value = 0.001036636855560411
$ws.Cells.Item(3,20).Value = 0.001036636855560411
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf18"
$ws.Cells.Item(4,3).Value = "Fgfr3"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.1601763333333333
$ws.Cells.Item(4,8).Value = 0.480529
$ws.Cells.Item(4,9).Value = 0.01412814675921196
$ws.Cells.Item(4,10).Value = 0.01412814675921196
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.445350666666667
$ws.Cells.Item(4,14).Value = 4.336052
$ws.Cells.Item(4,15).Value = 0.1861038648836906
$ws.Cells.Item(4,16).Value = 0.1861038648836906
$ws.Cells.Item(4,17).Value = 0.2315109701675556
$ws.Cells.Item(4,18).Value = 2.083598731508
$ws.Cells.Item(4,19).Value = 0.002629302715533334
$ws.Cells.Item(4,20).Value = 0.002629302715533334
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fgf18"
$ws.Cells.Item(5,3).Value = "Fgfr3"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 10.019353
$ws.Cells.Item(5,8).Value = 30.058059
$ws.Cells.Item(5,9).Value = 0.8837441004581448
$ws.Cells.Item(5,10).Value = 0.8837441004581448
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 5.751166666666666
$ws.Cells.Item(5,14).Value = 17.2535
$ws.Cells.Item(5,15).Value = 0.7405222614421495
$ws.Cells.Item(5,16).Value = 0.7405222614421495
$ws.Cells.Item(5,17).Value = 57.62296899516667
$ws.Cells.Item(5,18).Value = 518.6067209565
$ws.Cells.Item(5,19).Value = 0.6544321798074235
$ws.Cells.Item(5,20).Value = 0.6544321798074235
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf18"
$ws.Cells.Item(6,3).Value = "Fgfr3"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 10.019353
$ws.Cells.Item(6,8).Value = 30.058059
$ws.Cells.Item(6,9).Value = 0.8837441004581448
$ws.Cells.Item(6,10).Value = 0.8837441004581448
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.5698483333333333
$ws.Cells.Item(6,14).Value = 1.709545
$ws.Cells.Item(6,15).Value = 0.07337387367415998
$ws.Cells.Item(6,16).Value = 0.07337387367416
$ws.Cells.Item(6,17).Value = 5.709511608128333
$ws.Cells.Item(6,18).Value = 51.385604473155
$ws.Cells.Item(6,19).Value = 0.06484372798730006
$ws.Cells.Item(6,20).Value = 0.06484372798730008
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf18"
$ws.Cells.Item(7,3).Value = "Fgfr3"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 10.019353
$ws.Cells.Item(7,8).Value = 30.058059
$ws.Cells.Item(7,9).Value = 0.8837441004581448
$ws.Cells.Item(7,10).Value = 0.8837441004581448
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.445350666666667
$ws.Cells.Item(7,14).Value = 4.336052
$ws.Cells.Item(7,15).Value = 0.1861038648836906
$ws.Cells.Item(7,16).Value = 0.1861038648836906
$ws.Cells.Item(7,17).Value = 14.48147853811867
$ws.Cells.Item(7,18).Value = 130.333306843068
$ws.Cells.Item(7,19).Value = 0.1644681926634213
$ws.Cells.Item(7,20).Value = 0.1644681926634213
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Fgf18"
$ws.Cells.Item(8,3).Value = "Fgfr3"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 1.010195666666666
$ws.Cells.Item(8,8).Value = 3.030587
$ws.Cells.Item(8,9).Value = 0.08910300502687639
$ws.Cells.Item(8,10).Value = 0.0891030050268764
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 5.751166666666666
$ws.Cells.Item(8,14).Value = 17.2535
$ws.Cells.Item(8,15).Value = 0.7405222614421495
$ws.Cells.Item(8,16).Value = 0.7405222614421495
$ws.Cells.Item(8,17).Value = 5.809803644944443
$ws.Cells.Item(8,18).Value = 52.28823280449999
$ws.Cells.Item(8,19).Value = 0.06598275878379371
$ws.Cells.Item(8,20).Value = 0.06598275878379373
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Fgf18"
$ws.Cells.Item(9,3).Value = "Fgfr3"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 1.010195666666666
$ws.Cells.Item(9,8).Value = 3.030587
$ws.Cells.Item(9,9).Value = 0.08910300502687639
$ws.Cells.Item(9,10).Value = 0.0891030050268764
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.5698483333333333
$ws.Cells.Item(9,14).Value = 1.709545
$ws.Cells.Item(9,15).Value = 0.07337387367415998
$ws.Cells.Item(9,16).Value = 0.07337387367416
$ws.Cells.Item(9,17).Value = 0.5756583169905554
$ws.Cells.Item(9,18).Value = 5.180924852914999
$ws.Cells.Item(9,19).Value = 0.00653783263483007
$ws.Cells.Item(9,20).Value = 0.006537832634830073
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Fgf18"
$ws.Cells.Item(10,3).Value = "Fgfr3"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 1.010195666666666
$ws.Cells.Item(10,8).Value = 3.030587
$ws.Cells.Item(10,9).Value = 0.08910300502687639
$ws.Cells.Item(10,10).Value = 0.0891030050268764
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.445350666666667
$ws.Cells.Item(10,14).Value = 4.336052
$ws.Cells.Item(10,15).Value = 0.1861038648836906
$ws.Cells.Item(10,16).Value = 0.1861038648836906
$ws.Cells.Item(10,17).Value = 1.460086980280444
$ws.Cells.Item(10,18).Value = 13.140782822524
$ws.Cells.Item(10,19).Value = 0.01658241360825261
$ws.Cells.Item(10,20).Value = 0.01658241360825261
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Fgf18"
$ws.Cells.Item(11,3).Value = "Fgfr3"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.1476666666666667
$ws.Cells.Item(11,8).Value = 0.443
$ws.Cells.Item(11,9).Value = 0.01302474775576687
$ws.Cells.Item(11,10).Value = 0.01302474775576687
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 5.751166666666666
$ws.Cells.Item(11,14).Value = 17.2535
$ws.Cells.Item(11,15).Value = 0.7405222614421495
$ws.Cells.Item(11,16).Value = 0.7405222614421495
$ws.Cells.Item(11,17).Value = 0.849255611111111
$ws.Cells.Item(11,18).Value = 7.6433005
$ws.Cells.Item(11,19).Value = 0.009645115662814042
$ws.Cells.Item(11,20).Value = 0.009645115662814042
$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,2).Value = "Fgf18"
$ws.Cells.Item(12,3).Value = "Fgfr3"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.1476666666666667
$ws.Cells.Item(12,8).Value = 0.443
$ws.Cells.Item(12,9).Value = 0.01302474775576687
$ws.Cells.Item(12,10).Value = 0.01302474775576687
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.5698483333333333
$ws.Cells.Item(12,14).Value = 1.709545
$ws.Cells.Item(12,15).Value = 0.07337387367415998
$ws.Cells.Item(12,16).Value = 0.07337387367416
$ws.Cells.Item(12,17).Value = 0.08414760388888888
$ws.Cells.Item(12,18).Value = 0.7573284349999999
$ws.Cells.Item(12,19).Value = 0.000955676196469437
$ws.Cells.Item(12,20).Value = 0.0009556761964694372
$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,2).Value = "Fgf18"
$ws.Cells.Item(13,3).Value = "Fgfr3"
$ws.Cells.Item(13,4).Value = "MuSCs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.1476666666666667
$ws.Cells.Item(13,8).Value = 0.443
$ws.Cells.Item(13,9).Value = 0.01302474775576687
$ws.Cells.Item(13,10).Value = 0.01302474775576687
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.445350666666667
$ws.Cells.Item(13,14).Value = 4.336052
$ws.Cells.Item(13,15).Value = 0.1861038648836906
$ws.Cells.Item(13,16).Value = 0.1861038648836906
$ws.Cells.Item(13,17).Value = 0.2134301151111112
$ws.Cells.Item(13,18).Value = 1.920871036
$ws.Cells.Item(13,19).Value = 0.00242395589648339
$ws.Cells.Item(13,20).Value = 0.00242395589648339

Write-Output "Updated Fgf18-Fgfr3 LR-pair data with new TPM values"
